# Add "optimum targeting effect" slide (targeting.pptx commit)
$p = $ppt.ActivePresentation

# --- New slide 2, "Title and Content" layout (same layout index as slideLayout2.xml) ---
$s2 = $p.Slides.Add(2, 2)

$s2.Shapes.Item(1).TextFrame.TextRange.Text = "一律介入と最適ターゲティングの効果"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "アウトカム：返信時に提供意向を示したかどうか"

# --- Results table (6 rows x 3 cols) ---
$gf = $s2.Shapes.AddTable(6, 3, 72, 144, 720, 432)
$tbl = $gf.Table

# column widths (EMU / 12700 = points)
$tbl.Columns.Item(1).Width = 1828800 / 12700
$tbl.Columns.Item(2).Width = 914400 / 12700
$tbl.Columns.Item(3).Width = 914400 / 12700

# row heights
for ($r = 1; $r -le 6; $r++) {
    $tbl.Rows.Item($r).Height = 228600 / 12700
}

# turn off the default banded-table look so plain borders read cleanly
$tbl.FirstRow = $false
$tbl.HorizBanding = $false

$data = @(
    @("", "Predicted treatment effect", "Predicted treatment effect"),
    @(" ", "Mean", "SD"),
    @("Treatment B (uniform)", "0.0206", "0.0570"),
    @("Treatment C (uniform)", "-0.0019", "0.0598"),
    @("Treatment D (uniform)", "0.0079", "0.0628"),
    @("Optimum targeting", "0.0473", "0.0530")
)

for ($r = 1; $r -le 6; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $cell = $tbl.Cell($r, $c)
        $tr = $cell.Shape.TextFrame.TextRange
        $tr.Text = $data[$r - 1][$c - 1]
        $tr.Font.Name = "Arial"
        $tr.Font.Size = 11
        $tr.Font.Color.RGB = 0
        if ($c -eq 1) {
            $tr.ParagraphFormat.Alignment = 1
        } else {
            $tr.ParagraphFormat.Alignment = 2
        }

        # clear every border, then draw back only the ones the design calls for
        for ($b = 1; $b -le 4; $b++) {
            $cell.Borders.Item($b).Visible = $false
        }
    }
}

# merge the top-right header cells ("Predicted treatment effect")
$tbl.Cell(1, 2).Merge($tbl.Cell(1, 3))

# header rule above row 1, rule below row 1 (col 2-3) and below row 2 (full width), rule below row 6
for ($c = 1; $c -le 3; $c++) {
    $top = $tbl.Cell(1, $c).Borders.Item(1)
    $top.Visible = $true
    $top.Weight = 1
    $top.ForeColor.RGB = 0
}
for ($c = 2; $c -le 3; $c++) {
    $bot = $tbl.Cell(1, $c).Borders.Item(3)
    $bot.Visible = $true
    $bot.Weight = 1
    $bot.ForeColor.RGB = 0
}
for ($c = 1; $c -le 3; $c++) {
    $bot = $tbl.Cell(2, $c).Borders.Item(3)
    $bot.Visible = $true
    $bot.Weight = 1
    $bot.ForeColor.RGB = 0

    $bot6 = $tbl.Cell(6, $c).Borders.Item(3)
    $bot6.Visible = $true
    $bot6.Weight = 1
    $bot6.ForeColor.RGB = 0
}

Write-Host "optimum targeting slide added"
